# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de):
#   - mark the two files as handed back (Status: "Ready for handoff" ->
#     "Handed back: in sync with en-US", on the Overview sheet too, since it
#     shares the same text)
#   - add "Latest Target File" (E) / "Latest Handback File" (F) columns,
#     populated (+ hyperlinked) with the same file references already shown
#     under "Latest Handoff File" (C) and the source file (A)
#   - stamp "Latest Handback DateTime" (G) with the real handback time,
#     replacing the "0001-01-01 00:00:00" placeholder

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: just the status text (B/C columns, rows 2-3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- per-language handback datetimes ---
$handbackDateTime = @{
    "zh-cn" = "2016-03-09 16:04:21"
    "de-de" = "2016-03-09 16:04:35"
}

foreach ($sheetName in "zh-cn", "de-de") {
    $ws = $wb.Worksheets.Item($sheetName)

    # Snapshot the existing hyperlink addresses for row 2 / row 3 before we
    # start touching the sheet (A2, C2, A3, C3 in that collection order).
    $addrs = @()
    foreach ($h in $ws.Hyperlinks) {
        $addrs += $h.Address
    }
    $mdAddr2  = $addrs[0]
    $xlfAddr2 = $addrs[1]
    $mdAddr3  = $addrs[2]
    $xlfAddr3 = $addrs[3]

    $mdDisplay2  = $ws.Range("A2").Value2
    $xlfDisplay2 = $ws.Range("C2").Value2
    $mdDisplay3  = $ws.Range("A3").Value2
    $xlfDisplay3 = $ws.Range("C3").Value2

    # Status -> handed back
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # New "Latest Target File" / "Latest Handback File" columns, row 2
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdAddr2, $null, $null, $mdDisplay2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfAddr2, $null, $null, $xlfDisplay2)

    # New "Latest Target File" / "Latest Handback File" columns, row 3
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdAddr3, $null, $null, $mdDisplay3)
    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfAddr3, $null, $null, $xlfDisplay3)

    # Latest Handback DateTime, rows 2 & 3
    $ws.Range("G2").Value = $handbackDateTime[$sheetName]
    $ws.Range("G3").Value = $handbackDateTime[$sheetName]
}
